$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86, shifting the existing rows 86..119 down to 87..120
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new weekly data point
$ws.Range("A86").Value = 10
$ws.Range("B86").Value = "Vega Modelo de Temuco"
$ws.Range("C86").Value = "La Araucanía"
$ws.Range("D86").Value = 44523
$ws.Range("E86").Value = 9
$ws.Range("F86").Value = 100114007
$ws.Range("G86").Value = "Jengibre"
$ws.Range("H86").Value = "Sin especificar"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 50
$ws.Range("K86").Value = 20000
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = 20000
$ws.Range("N86").Value = "$/caja 13 kilos"
$ws.Range("O86").Value = "Perú"
$ws.Range("P86").Value = 1538
$ws.Range("Q86").Value = 13
$ws.Range("R86").Value = "Hortaliza"
